# Update "Förändrad" (Changed) date column (C) for rows 2-8 from 45183 (2023-09-14)
# to 45184 (2023-09-15) to reflect the automatic daily refresh of this logging report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45184   # Column C
}
